$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 443.33334
$ws.Cells.Item(17, 10).Value = 443.33334
$ws.Cells.Item(17, 12).Value = 1330.00002
$ws.Cells.Item(17, 14).Value = -1666.00002
$ws.Cells.Item(18, 8).Value = 288.29413
$ws.Cells.Item(18, 9).Value = 292.3846
$ws.Cells.Item(18, 10).Value = 275
$ws.Cells.Item(18, 11).Value = 292.3846
$ws.Cells.Item(18, 12).Value = 275
$ws.Cells.Item(18, 13).Value = -8.384599999999978
$ws.Cells.Item(18, 14).Value = -843
$ws.Cells.Item(137, 8).Value = 1321.1818
$ws.Cells.Item(137, 9).Value = 1140.9412
$ws.Cells.Item(137, 10).Value = 1512.6875
$ws.Cells.Item(137, 11).Value = 3422.8236
$ws.Cells.Item(137, 12).Value = 4538.0625
$ws.Cells.Item(137, 13).Value = -872.8235999999997
$ws.Cells.Item(137, 14).Value = -9638.0625
$ws.Cells.Item(141, 8).Value = 8777.888999999999
$ws.Cells.Item(141, 9).Value = 4698.6665
$ws.Cells.Item(141, 10).Value = 10817.5
$ws.Cells.Item(141, 11).Value = 14095.9995
$ws.Cells.Item(141, 12).Value = 32452.5
$ws.Cells.Item(141, 13).Value = -8915.999500000002
$ws.Cells.Item(141, 14).Value = -42812.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21282060
$ws.Cells.Item(32, 9).Value = 5250.8374
$ws.Cells.Item(32, 10).Value = 250007740
$ws.Cells.Item(32, 11).Value = 5250.8374
$ws.Cells.Item(32, 12).Value = 250007740
$ws.Cells.Item(32, 13).Value = -4963.8374
$ws.Cells.Item(32, 14).Value = -250008314
$ws.Cells.Item(61, 8).Value = 11112459
$ws.Cells.Item(61, 9).Value = 11112459
$ws.Cells.Item(61, 11).Value = 11112459
$ws.Cells.Item(61, 13).Value = -11112247
$ws.Cells.Item(74, 8).Value = 963.59576
$ws.Cells.Item(74, 9).Value = 855
$ws.Cells.Item(74, 10).Value = 1422.1111
$ws.Cells.Item(74, 11).Value = 855
$ws.Cells.Item(74, 12).Value = 1422.1111
$ws.Cells.Item(74, 13).Value = 19
$ws.Cells.Item(74, 14).Value = -3170.1111
$ws.Cells.Item(77, 8).Value = 963.59576
$ws.Cells.Item(77, 9).Value = 855
$ws.Cells.Item(77, 10).Value = 1422.1111
$ws.Cells.Item(77, 11).Value = 4275
$ws.Cells.Item(77, 12).Value = 7110.5555
$ws.Cells.Item(77, 13).Value = 93
$ws.Cells.Item(77, 14).Value = -15846.5555
$ws.Cells.Item(97, 8).Value = 941.55554
$ws.Cells.Item(97, 9).Value = 921.125
$ws.Cells.Item(97, 10).Value = 1105
$ws.Cells.Item(97, 11).Value = 921.125
$ws.Cells.Item(97, 12).Value = 1105
$ws.Cells.Item(97, 13).Value = -425.125
$ws.Cells.Item(97, 14).Value = -2097
$ws.Cells.Item(132, 8).Value = 3679440.8
$ws.Cells.Item(132, 9).Value = 3055.6365
$ws.Cells.Item(132, 10).Value = 11767488
$ws.Cells.Item(132, 11).Value = 9166.9095
$ws.Cells.Item(132, 12).Value = 35302464
$ws.Cells.Item(132, 13).Value = -6636.9095
$ws.Cells.Item(132, 14).Value = -35307524
$ws.Cells.Item(136, 8).Value = 11112459
$ws.Cells.Item(136, 9).Value = 11112459
$ws.Cells.Item(136, 11).Value = 33337377
$ws.Cells.Item(136, 13).Value = -33334827
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1293548.9
$ws.Cells.Item(86, 9).Value = 1461.8
$ws.Cells.Item(86, 10).Value = 2908657.8
$ws.Cells.Item(86, 11).Value = 1461.8
$ws.Cells.Item(86, 12).Value = 2908657.8
$ws.Cells.Item(86, 13).Value = -338.8
$ws.Cells.Item(86, 14).Value = -2910903.8
$ws.Cells.Item(89, 8).Value = 1293548.9
$ws.Cells.Item(89, 9).Value = 1461.8
$ws.Cells.Item(89, 10).Value = 2908657.8
$ws.Cells.Item(89, 11).Value = 7309
$ws.Cells.Item(89, 12).Value = 14543289
$ws.Cells.Item(89, 13).Value = -1693
$ws.Cells.Item(89, 14).Value = -14554521
$ws.Cells.Item(94, 8).Value = 648.2083
$ws.Cells.Item(94, 9).Value = 611.85
$ws.Cells.Item(94, 10).Value = 830
$ws.Cells.Item(94, 11).Value = 611.85
$ws.Cells.Item(94, 12).Value = 830
$ws.Cells.Item(94, 13).Value = -160.85
$ws.Cells.Item(94, 14).Value = -1732
$ws.Cells.Item(134, 8).Value = 4189.75
$ws.Cells.Item(134, 9).Value = 1071.1482
$ws.Cells.Item(134, 11).Value = 3213.4446
$ws.Cells.Item(134, 13).Value = -678.4446000000003
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1134788.6
$ws.Cells.Item(31, 9).Value = 1462848.6
$ws.Cells.Item(31, 10).Value = 1490.8182
$ws.Cells.Item(31, 11).Value = 1462848.6
$ws.Cells.Item(31, 12).Value = 1490.8182
$ws.Cells.Item(31, 13).Value = -1462553.6
$ws.Cells.Item(31, 14).Value = -2080.8182
$ws.Cells.Item(34, 8).Value = 1134788.6
$ws.Cells.Item(34, 9).Value = 1462848.6
$ws.Cells.Item(34, 10).Value = 1490.8182
$ws.Cells.Item(34, 11).Value = 1462848.6
$ws.Cells.Item(34, 12).Value = 1490.8182
$ws.Cells.Item(34, 13).Value = -1462646.6
$ws.Cells.Item(34, 14).Value = -1894.8182
$ws.Cells.Item(58, 8).Value = 31250670
$ws.Cells.Item(58, 9).Value = 41667292
$ws.Cells.Item(58, 10).Value = 803.5
$ws.Cells.Item(58, 11).Value = 41667292
$ws.Cells.Item(58, 12).Value = 803.5
$ws.Cells.Item(58, 13).Value = -41667089
$ws.Cells.Item(58, 14).Value = -1209.5
$ws.Cells.Item(134, 8).Value = 1337.7
$ws.Cells.Item(134, 9).Value = 1248.1538
$ws.Cells.Item(134, 11).Value = 3744.4614
$ws.Cells.Item(134, 13).Value = -1209.4614
$ws.Cells.Item(136, 8).Value = 31250670
$ws.Cells.Item(136, 9).Value = 41667292
$ws.Cells.Item(136, 10).Value = 803.5
$ws.Cells.Item(136, 11).Value = 125001876
$ws.Cells.Item(136, 12).Value = 2410.5
$ws.Cells.Item(136, 13).Value = -124999326
$ws.Cells.Item(136, 14).Value = -7510.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 1356.8125
$ws.Cells.Item(129, 9).Value = 920.9
$ws.Cells.Item(129, 10).Value = 2083.3333
$ws.Cells.Item(129, 11).Value = 2762.7
$ws.Cells.Item(129, 12).Value = 6249.999899999999
$ws.Cells.Item(129, 13).Value = 2237.3
$ws.Cells.Item(129, 14).Value = -16249.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 8270.579
$ws.Cells.Item(132, 9).Value = 2912.9
$ws.Cells.Item(132, 10).Value = 14223.556
$ws.Cells.Item(132, 11).Value = 8738.700000000001
$ws.Cells.Item(132, 12).Value = 42670.66800000001
$ws.Cells.Item(132, 13).Value = -6208.700000000001
$ws.Cells.Item(132, 14).Value = -47730.66800000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 4
$ws.Cells.Item(14, 9).Value = 4
$ws.Cells.Item(14, 11).Value = 4
$ws.Cells.Item(14, 13).Value = 168
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(40, 8).Value = 20834658
$ws.Cells.Item(40, 9).Value = 1344.3334
$ws.Cells.Item(40, 10).Value = 83334600
$ws.Cells.Item(40, 11).Value = 1344.3334
$ws.Cells.Item(40, 12).Value = 83334600
$ws.Cells.Item(40, 13).Value = -1208.3334
$ws.Cells.Item(40, 14).Value = -83334872
$ws.Cells.Item(54, 8).Value = 7250
$ws.Cells.Item(54, 10).Value = 7250
$ws.Cells.Item(54, 12).Value = 7250
$ws.Cells.Item(54, 14).Value = -8538
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 25281.428
$ws.Cells.Item(12, 9).Value = 1992
$ws.Cells.Item(12, 10).Value = 42748.5
$ws.Cells.Item(12, 11).Value = 1992
$ws.Cells.Item(12, 12).Value = 42748.5
$ws.Cells.Item(12, 13).Value = -1850
$ws.Cells.Item(12, 14).Value = -43032.5
$ws.Cells.Item(15, 8).Value = 9800
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 9800
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 13).Value = 9800
$ws.Cells.Item(15, 14).Value = -10376
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(24, 8).Value = 3000
$ws.Cells.Item(24, 10).Value = 3000
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 14).Value = -3460
$ws.Cells.Item(29, 8).Value = 3792
$ws.Cells.Item(29, 9).Value = 3000
$ws.Cells.Item(29, 10).Value = 3990
$ws.Cells.Item(29, 11).Value = 3000
$ws.Cells.Item(29, 12).Value = 3990
$ws.Cells.Item(29, 13).Value = -2710
$ws.Cells.Item(29, 14).Value = -4570
$ws.Cells.Item(51, 8).Value = 4100
$ws.Cells.Item(51, 10).Value = 4100
$ws.Cells.Item(51, 12).Value = 4100
$ws.Cells.Item(51, 14).Value = -5120
$ws.Cells.Item(132, 8).Value = 73950.625
$ws.Cells.Item(132, 9).Value = 150744.58
$ws.Cells.Item(132, 10).Value = 14222
$ws.Cells.Item(132, 11).Value = 452233.74
$ws.Cells.Item(132, 12).Value = 42666
$ws.Cells.Item(132, 13).Value = -449703.74
$ws.Cells.Item(132, 14).Value = -47726
$ws.Cells.Item(136, 8).Value = 2040.7106
$ws.Cells.Item(136, 9).Value = 1359.3158
$ws.Cells.Item(136, 10).Value = 2722.1052
$ws.Cells.Item(136, 11).Value = 4077.9474
$ws.Cells.Item(136, 12).Value = 8166.3156
$ws.Cells.Item(136, 13).Value = -1527.9474
$ws.Cells.Item(136, 14).Value = -13266.3156
